# aggiornamento a 9/09 compreso
# Appends 8 new daily rows (2021-09-02 .. 2021-09-09) to the bottom of the
# report, matching the column A date-serial / B "nuovi pos." / C "somma
# mobile 7gg." / D "somma mobile 7gg. per 100mila abitanti" layout already
# used by the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$data = @(
  @(44441, 4, 15, 85.37765382207297),
  @(44442, 3, 14, 79.68581023393477),
  @(44443, 2, 14, 79.68581023393477),
  @(44444, 2, 16, 91.06949741021117),
  @(44445, 2, 13, 73.99396664579658),
  @(44446, 0, 13, 73.99396664579658),
  @(44447, 0, 13, 73.99396664579658),
  @(44448, 0, 9, 51.22659229324378)
)

$lastRow = 366
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i

  # Column A carries the same date formatting/border/alignment style as the
  # rest of the column; copy it down from the previous row before writing
  # the new value so the new cell picks up style index 2 instead of Excel's
  # default style.
  $ws.Range("A$lastRow").Copy()
  $ws.Range("A$row").PasteSpecial(-4122)

  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
  $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$excel.CutCopyMode = 0
